$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "65.117.64"
$ws.Cells.Item(2, 5).Value = "  +0.09%  "
$ws.Cells.Item(3, 4).Value = "3.425.86"
$ws.Cells.Item(3, 5).Value = "  +1.83%  "
$ws.Cells.Item(4, 5).Value = "  +0.42%  "
$ws.Cells.Item(5, 4).Value = "'548.39"
$ws.Cells.Item(5, 5).Value = "  +2.21%  "
$ws.Cells.Item(6, 4).Value = "'178.54"
$ws.Cells.Item(6, 5).Value = "  -1.54%  "
$ws.Cells.Item(7, 4).Value = "'0.632"
$ws.Cells.Item(7, 5).Value = "  +5.66%  "
$ws.Cells.Item(8, 5).Value = "  +0.06%  "
$ws.Cells.Item(9, 4).Value = "'0.625"
$ws.Cells.Item(9, 5).Value = "  +1.25%  "
$ws.Cells.Item(10, 4).Value = "'0.150"
$ws.Cells.Item(10, 5).Value = "  +7.62%  "
$ws.Cells.Item(11, 4).Value = "'53.37"
$ws.Cells.Item(11, 5).Value = "  -3.50%  "
$ws.Cells.Item(12, 4).Value = "'0.0000269"
$ws.Cells.Item(12, 5).Value = "  +1.76%  "
$ws.Cells.Item(13, 4).Value = "'9.13"
$ws.Cells.Item(13, 5).Value = "  +0.24%  "
$ws.Cells.Item(14, 4).Value = "3.967.91"
$ws.Cells.Item(14, 5).Value = "  +3.04%  "
$ws.Cells.Item(15, 5).Value = "  +1.40%  "
$ws.Cells.Item(16, 4).Value = "3.419.41"
$ws.Cells.Item(16, 5).Value = "  +2.83%  "
$ws.Cells.Item(17, 4).Value = "'18.23"
$ws.Cells.Item(17, 5).Value = "  +2.09%  "
$ws.Cells.Item(18, 4).Value = "65.186.09"
$ws.Cells.Item(18, 5).Value = "  +0.30%  "
$ws.Cells.Item(19, 4).Value = "'11.77"
$ws.Cells.Item(20, 4).Value = "'0.978"
$ws.Cells.Item(20, 5).Value = "  +0.15%  "
$ws.Cells.Item(21, 4).Value = "'412.69"
$ws.Cells.Item(21, 5).Value = "  +6.91%  "
$ws.Cells.Item(22, 4).Value = "'3.99"
$ws.Cells.Item(22, 5).Value = "  +6.04%  "
$ws.Cells.Item(23, 4).Value = "'4.26"
$ws.Cells.Item(23, 5).Value = "  +1.05%  "
$ws.Cells.Item(24, 4).Value = "'84.52"
$ws.Cells.Item(24, 5).Value = "  +2.28%  "
$ws.Cells.Item(25, 4).Value = "'10.76"
$ws.Cells.Item(25, 5).Value = "  -7.69%  "
$ws.Cells.Item(26, 4).Value = "'2.84"
$ws.Cells.Item(26, 5).Value = "  +1.67%  "
$ws.Cells.Item(27, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(27, 4).Value = "'12.14"
$ws.Cells.Item(27, 5).Value = "  +5.56%  "
$ws.Cells.Item(28, 2).Value = "LEO"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Cells.Item(28, 4).Value = "'6.01"
$ws.Cells.Item(28, 5).Value = "  -1.58%  "
$ws.Cells.Item(29, 4).Value = "'8.83"
$ws.Cells.Item(29, 5).Value = "  +5.31%  "
$ws.Cells.Item(30, 4).Value = "'29.67"
$ws.Cells.Item(30, 5).Value = "  +1.28%  "
$ws.Cells.Item(31, 2).Value = "NEARProtocol"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(31, 4).Value = "'6.51"
$ws.Cells.Item(31, 5).Value = "  -4.02%  "
$ws.Cells.Item(32, 2).Value = "Bittensor"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Cells.Item(32, 4).Value = "'610.27"
$ws.Cells.Item(32, 5).Value = "  -5.96%  "
$ws.Cells.Item(33, 4).Value = "'11.62"
$ws.Cells.Item(33, 5).Value = "  +2.30%  "
$ws.Cells.Item(34, 5).Value = "  +0.94%  "
$ws.Cells.Item(35, 4).Value = "'58.62"
$ws.Cells.Item(35, 5).Value = "  +1.44%  "
$ws.Cells.Item(36, 2).Value = "Kaspa"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(36, 4).Value = "'0.146"
$ws.Cells.Item(36, 5).Value = "  +14.20%  "
$ws.Cells.Item(37, 2).Value = "Dai"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(37, 4).Value = "'0.999"
$ws.Cells.Item(37, 5).Value = "  -0.28%  "
$ws.Cells.Item(38, 4).Value = "'37.13"
$ws.Cells.Item(38, 5).Value = "  -0.59%  "
$ws.Cells.Item(39, 4).Value = "0.0₃0771"
$ws.Cells.Item(39, 5).Value = "  -1.07%  "
$ws.Cells.Item(40, 4).Value = "'0.376"
$ws.Cells.Item(40, 5).Value = "  -3.77%  "
$ws.Cells.Item(41, 4).Value = "3.172.89"
$ws.Cells.Item(41, 5).Value = "  +5.60%  "
$ws.Cells.Item(42, 4).Value = "'3.30"
$ws.Cells.Item(42, 5).Value = "  +1.03%  "
$ws.Cells.Item(43, 4).Value = "'0.998"
$ws.Cells.Item(43, 5).Value = "  +0.30%  "
$ws.Cells.Item(44, 4).Value = "'2.53"
$ws.Cells.Item(44, 5).Value = "  -6.93%  "
$ws.Cells.Item(45, 5).Value = "  +2.04%  "
$ws.Cells.Item(46, 2).Value = "ApeXProtocol"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Cells.Item(46, 4).Value = "'3.22"
$ws.Cells.Item(46, 5).Value = "  +0.52%  "
$ws.Cells.Item(47, 4).Value = "'0.0408"
$ws.Cells.Item(47, 5).Value = "  -0.05%  "
$ws.Cells.Item(48, 2).Value = "WEMIXToken"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(48, 4).Value = "'2.71"
$ws.Cells.Item(48, 5).Value = "  +0.73%  "
$ws.Cells.Item(49, 4).Value = "'0.131"
$ws.Cells.Item(49, 5).Value = "  +3.35%  "
$ws.Cells.Item(50, 4).Value = "'137.86"
$ws.Cells.Item(50, 5).Value = "  -0.78%  "
$ws.Cells.Item(51, 4).Value = "'8.33"
$ws.Cells.Item(51, 5).Value = "  -0.98%  "
